$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First Integrated")

# Clear SWL Note (column H) for rows 2-94 (column H removed entirely)
for ($r = 2; $r -le 94; $r++) {
    $ws.Cells.Item($r, 8).Value = ""
}

# Set Manufacturer (column I) for specific rows
for ($r = 6; $r -le 12; $r++) {
    $ws.Cells.Item($r, 9).Value = "Tiger"
}
for ($r = 39; $r -le 45; $r++) {
    $ws.Cells.Item($r, 9).Value = "Ridge Gear"
}
for ($r = 48; $r -le 49; $r++) {
    $ws.Cells.Item($r, 9).Value = "Tiger"
}
for ($r = 54; $r -le 61; $r++) {
    $ws.Cells.Item($r, 9).Value = "Tiger"
}
$ws.Cells.Item(85, 9).Value = "Miller"
$ws.Cells.Item(91, 9).Value = "Ridge Gear"
$ws.Cells.Item(93, 9).Value = "Riley "
for ($r = 336; $r -le 417; $r++) {
    $ws.Cells.Item($r, 9).Value = "Crosby"
}
for ($r = 565; $r -le 642; $r++) {
    $ws.Cells.Item($r, 9).Value = "Crosby"
}

# Clear Model (column E) for rows 330-335 and 643-646 (A344 removed)
foreach ($r in 330..335) {
    $ws.Cells.Item($r, 5).Value = ""
}
foreach ($r in 643..646) {
    $ws.Cells.Item($r, 5).Value = ""
}

# Update Model (column E) for row 85: RHINO -> BLACK RHINO
$ws.Cells.Item(85, 5).Value = "BLACK RHINO"
